$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.552.95"
$ws.Range("E2").Value = "  -1.80%  "

$ws.Range("D3").Value = "2.627.18"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "2.635.69"
$ws.Range("E9").Value = "  +0.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("D14").Value = "3.097.72"
$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").Value = "58.507.97"
$ws.Range("E15").Value = "  -1.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("D17").Value = "2.655.71"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000132"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.75%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.60%  "

$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.81%  "

$ws.Range("D29").Value = "0.0₃0738"
$ws.Range("E29").Value = "  -1.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").Value = "  -1.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.21%  "

$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.54%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.840"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("E39").Value = "  -3.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "281.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.33%  "

$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.56%  "

$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0935"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0224"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").Value = "1.945.72"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
